$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table runs A1:R129, one header row + 64 weekly pairs of rows
# (Primera/Segunda quality) sorted newest-first. A new week's pair of
# readings was recorded, so insert two fresh rows right above the first
# data row (row 46) and let Excel push the rest of the table down.
$ws.Rows.Item(46).Insert()
$ws.Rows.Item(47).Insert()

# Seed the two new rows with the same record shape as the row that is
# now directly below them (the former row 46/47, now 48/49), then fix
# up the date for the new week.
$src = $ws.Range("A48:R49")
$dst = $ws.Range("A46:R47")
$src.Copy($dst)

$ws.Cells.Item(46, 4).Value2 = 44467
$ws.Cells.Item(47, 4).Value2 = 44467
